$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.480.46"
$ws.Range("E2").Value = "  -4.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.967.37"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.44"
$ws.Range("E5").Value = "  -5.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.69"
$ws.Range("E6").Value = "  -8.05%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.978.86"
$ws.Range("E9").Value = "  -5.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.112"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("E11").Value = "  -6.59%  "
$ws.Range("E12").Value = "  -4.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.490.89"
$ws.Range("E13").Value = "  -5.29%  "
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.585.83"
$ws.Range("E15").Value = "  -4.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.57"
$ws.Range("E16").Value = "  -6.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.977.10"
$ws.Range("E17").Value = "  -5.33%  "
$ws.Range("E18").Value = "  -4.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.15"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.99"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "379.85"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.65"
$ws.Range("E22").Value = "  -5.90%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.66"
$ws.Range("E24").Value = "  -2.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.38"
$ws.Range("E25").Value = "  -3.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.468"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.095.28"
$ws.Range("E27").Value = "  -5.36%  "
$ws.Range("E28").Value = "  -3.39%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0937"
$ws.Range("E30").Value = "  -6.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.20"
$ws.Range("E31").Value = "  -6.42%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  -5.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.39"
$ws.Range("E34").Value = "  -3.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.91"
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.63"
$ws.Range("E36").Value = "  -3.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.90"
$ws.Range("E37").Value = "  -5.65%  "
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("E39").Value = "  -5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.54"
$ws.Range("E40").Value = "  -7.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.49"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.403.79"
$ws.Range("E43").Value = "  -9.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.07"
$ws.Range("E44").Value = "  -6.13%  "
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0588"
$ws.Range("E46").Value = "  -3.55%  "
$ws.Range("E47").Value = "  -6.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  +0.18%  "
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0949"
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.64"
$ws.Range("E51").Value = "  -6.44%  "
